# Updated cryptos list on Thu Sep 14 06:37:08 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-formatted numbers (e.g. "212.03", "1.619.55")
# that must stay text, not be auto-coerced to numeric cells.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price / Volume(1h) updates for rows 2-45 (column D = Price, column E = Volume(1h)) ---
# Each entry: row, Price (or $null to leave unchanged), Volume (or $null to leave unchanged)
$updates = @(
    @(2,  "26.300.69",  "  +1.29%  "),
    @(3,  "1.619.55",   $null),
    @(4,  $null,        "  -0.01%  "),
    @(5,  "212.03",     "  +0.81%  "),
    @(6,  $null,        "  -0.04%  "),
    @(7,  $null,        "  +0.93%  "),
    @(8,  $null,        "  +0.68%  "),
    @(9,  $null,        "  +0.83%  "),
    @(10, "18.78",      "  +4.89%  "),
    @(11, "0.0815",     "  +0.96%  "),
    @(12, "1.844.78",   "  +1.96%  "),
    @(13, "1.614.08",   "  +1.62%  "),
    @(14, "4.01",       "  +0.79%  "),
    @(15, $null,        "  +1.55%  "),
    @(16, "26.306.52",  "  +1.36%  "),
    @(17, "62.22",      "  +3.67%  "),
    @(18, $null,        "  +1.10%  "),
    @(19, $null,        "  -0.04%  "),
    @(20, "201.45",     "  +1.06%  "),
    @(21, $null,        "  +1.68%  "),
    @(22, "9.33",       "  +1.74%  "),
    @(23, $null,        "  +1.29%  "),
    @(24, $null,        "  +2.49%  "),
    @(25, "144.66",     "  +1.52%  "),
    @(26, $null,        "  -0.04%  "),
    @(27, $null,        "  -1.40%  "),
    @(28, $null,        "  +0.69%  "),
    @(29, "6.56",       "  +1.75%  "),
    @(30, $null,        "  +10.05%  "),
    @(31, $null,        "  +0.83%  "),
    @(32, "3.18",       "  +2.02%  "),
    @(33, "2.92",       "  -0.22%  "),
    @(34, "1.49",       "  +1.94%  "),
    @(35, $null,        "  +2.11%  "),
    @(36, "1.179.91",   "  +5.16%  "),
    @(37, $null,        "  +0.54%  "),
    @(38, "0.805",      "  +2.98%  "),
    @(39, $null,        "  +0.00%  "),
    @(40, $null,        "  +0.36%  "),
    @(41, $null,        "  +1.61%  "),
    @(42, "0.788",      "  +1.42%  "),
    @(43, $null,        "  +4.91%  "),
    @(44, "1.756.09",   "  +2.07%  "),
    @(45, "92.68",      "  +0.86%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $price = $u[1]
    $volume = $u[2]
    if ($null -ne $price) {
        $ws.Cells.Item($row, 4).Value = $price
    }
    if ($null -ne $volume) {
        $ws.Cells.Item($row, 5).Value = $volume
    }
}

# --- New coin (BabyDogeCoin) inserted at row 46, shifting existing rows 46-51 down by one ---
# This pushes "EnergySwap" (previously row 51) out of the sheet entirely.
$babyDogePrice = "0.0₆0106"

$newRows = @(
    @(46, "BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", $babyDogePrice, "  +15.52%  "),
    @(47, "RenderToken",  "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr",  "1.54",          "  +3.54%  "),
    @(48, "Aave",         "https://coinranking.com/coin/ixgUfzmLR+aave-aave",               "53.72",         "  +0.95%  "),
    @(49, "Cronos",       "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro",           "0.0508",        "  +1.10%  "),
    @(50, "Mantle",       "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt",              "0.408",         "  +0.43%  "),
    @(51, "USDD",         "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd",               "1.00",          "  -0.21%  ")
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}

$wb.Save()
